# Update cryptos list (crypto prices / 1h volume changes) - GitHub Actions refresh.
# Cells whose new value could be mis-parsed by Excel as a pure number (which would
# silently change the cell's type from Text to Number and reformat the display,
# e.g. "1.00" -> 1) are first forced to Text format ("@") before the value is written,
# matching the original inline-string / text cells in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = '27.368.23'
$ws.Range("E2").Value = '  -1.67%  '

# Row 3 - Ethereum
$ws.Range("D3").Value = '1.656.12'
$ws.Range("E3").Value = '  -0.40%  '

# Row 4 - TetherUSD
$ws.Range("E4").Value = '  -0.17%  '

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.33'
$ws.Range("E5").Value = '  -0.63%  '

# Row 6 - XRP
$ws.Range("E6").Value = '  -0.19%  '

# Row 7 - USDC
$ws.Range("E7").Value = '  -0.18%  '

# Row 8 - Solana
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.61'
$ws.Range("E8").Value = '  +0.73%  '

# Row 9 - Cardano
$ws.Range("E9").Value = '  -0.24%  '

# Row 10 - Dogecoin
$ws.Range("E10").Value = '  -1.06%  '

# Row 11 - TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0875'
$ws.Range("E11").Value = '  -0.44%  '

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = '1.890.11'
$ws.Range("E12").Value = '  -0.54%  '

# Row 13 - WrappedEther
$ws.Range("D13").Value = '1.655.38'
$ws.Range("E13").Value = '  -0.84%  '

# Row 14 - was Polygon, now Polkadot (rows 14/15 swapped order/ranking)
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.07'
$ws.Range("E14").Value = '  -1.52%  '

# Row 15 - was Polkadot, now Polygon
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.569'
$ws.Range("E15").Value = '  +3.80%  '

# Row 16 - Litecoin
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.71'
$ws.Range("E16").Value = '  -0.39%  '

# Row 17 - WrappedBTC
$ws.Range("D17").Value = '27.382.43'
$ws.Range("E17").Value = '  -1.53%  '

# Row 18 - BitcoinCash
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '231.85'
$ws.Range("E18").Value = '  -6.89%  '

# Row 19 - ShibaInu
$ws.Range("D19").Value = '0.0₃0726'
$ws.Range("E19").Value = '  -0.53%  '

# Row 20 - Chainlink
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.49'
$ws.Range("E20").Value = '  -0.05%  '

# Row 21 - Dai
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.00'
$ws.Range("E21").Value = '  -0.08%  '

# Row 22 - Uniswap
$ws.Range("E22").Value = '  -2.01%  '

# Row 23 - Avalanche
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.37'
$ws.Range("E23").Value = '  +0.47%  '

# Row 24 - Toncoin
$ws.Range("E24").Value = '  -1.04%  '

# Row 25 - Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.60'
$ws.Range("E25").Value = '  +0.67%  '

# Row 26 - Cosmos
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.12'
$ws.Range("E26").Value = '  -0.93%  '

# Row 27 - EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.91'
$ws.Range("E27").Value = '  -2.17%  '

# Row 28 - BinanceUSD
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  -0.20%  '

# Row 29 - Stellar
$ws.Range("E29").Value = '  -0.23%  '

# Row 30 - Hedera
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0496'
$ws.Range("E30").Value = '  -0.59%  '

# Row 31 - PancakeSwap
$ws.Range("E31").Value = '  -4.19%  '

# Row 32 - Filecoin
$ws.Range("E32").Value = '  -1.24%  '

# Row 33 - Maker
$ws.Range("D33").Value = '1.438.55'
$ws.Range("E33").Value = '  +0.42%  '

# Row 34 - InternetComputer(DFINITY)
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.15'
$ws.Range("E34").Value = '  +0.70%  '

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = '  +0.66%  '

# Row 36 - HuobiToken
$ws.Range("E36").Value = '  -0.71%  '

# Row 37 - ARBITRUM
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.909'
$ws.Range("E37").Value = '  -2.03%  '

# Row 38 - ImmutableX
$ws.Range("E38").Value = '  -1.89%  '

# Row 39 - VeChain
$ws.Range("E39").Value = '  +0.36%  '

# Row 40 - WEMIXToken
$ws.Range("E40").Value = '  -0.12%  '

# Row 41 - PaxDollar
$ws.Range("E41").Value = '  -0.19%  '

# Row 42 - FraxShare
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.52'
$ws.Range("E42").Value = '  +2.10%  '

# Row 43 - Aave
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '65.05'
$ws.Range("E43").Value = '  -6.67%  '

# Row 44 - MXToken
$ws.Range("E44").Value = '  +0.03%  '

# Row 45 - TrustWalletToken
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.789'
$ws.Range("E45").Value = '  +0.28%  '

# Row 46 - RocketPoolETH
$ws.Range("D46").Value = '1.798.04'
$ws.Range("E46").Value = '  -0.54%  '

# Row 47 - RenderToken
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.69'
$ws.Range("E47").Value = '  -0.72%  '

# Row 48 - Quant
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '88.04'
$ws.Range("E48").Value = '  -1.25%  '

# Row 49 - BabyDogeCoin
$ws.Range("E49").Value = '  -2.68%  '

# Row 50 - Algorand
$ws.Range("E50").Value = '  -0.21%  '

# Row 51 - EnergySwap
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.72'
$ws.Range("E51").Value = '  -1.11%  '
